# Update xpath selectors in the "menu" and "body" sheets to use the
# new //*[@id="root"]/... and //*[@id="footer"] style xpaths, and test
# the remaining selections/suitcases by moving the active selection.

$wb = $excel.ActiveWorkbook

$menu = $wb.Worksheets.Item("menu")
$body = $wb.Worksheets.Item("body")

# menu sheet: logo xpath (C2) and section_login xpath (C7)
$menu.Range("C2").Value = '//*[@id="root"]/div/div/div[1]/div[1]/div/nav/div[1]/a/img'
$menu.Range("C7").Value = '//*[@id="root"]/div/div/div[1]/div[1]/div/nav/div[3]/button'

# body sheet: logo xpath (C3) and footer xpath (C7)
$body.Range("C3").Value = '//*[@id="root"]/div/div/div[1]/div[1]/div/nav/div[1]/a/img'
$body.Range("C7").Value = '//*[@id="footer"]'

# Update the active selections on each sheet as recorded when the
# workbook was saved after testing.
$body.Activate()
$body.Range("C8").Select()

$menu.Activate()
$menu.Range("C2").Select()

$body.Activate()
